$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend "stim/" to each .bmp filename in column C (rows 2-9)
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = "stim/" + $cell.Value2
}

# Update the selected cell from E12 to C9
$ws.Range("C9").Select()
